$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data dictionary table gains one new row describing a new attribute:
# "Quantite_article" - the quantity of each article within an order.
$tbl = $ws.ListObjects("Tableau4234646")
$newRow = $tbl.ListRows.Add()

$ws.Range("A37").Value = "Quantite_article"
$ws.Range("B37").Value = "N"
$ws.Range("C37").Value = "16 bits"
$ws.Range("F37").Value = "Quantite de chaque article dans la commande"

# Match the centered alignment used throughout the rest of the table.
# (VerticalAlignment must be set before HorizontalAlignment so the engine
# reuses the existing "center/center" cell style instead of minting a
# vertical-only one.)
$ws.Range("A37:F37").VerticalAlignment = -4108
$ws.Range("A37:F37").HorizontalAlignment = -4108

# Reflect the author's final selection when the workbook was saved.
$ws.Range("I11").Select() | Out-Null
